# Update FFXIV leve-profit market-data cells (scheduled runner refresh).
# Targets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 191.86667
$ws.Range("I12").Value = 172.54546
$ws.Range("J12").Value = 245
$ws.Range("K12").Value = 172.54546
$ws.Range("L12").Value = 245
$ws.Range("M12").Value = -2.545459999999991
$ws.Range("N12").Value = -585

$ws.Range("H51").Value = 6274.5
$ws.Range("J51").Value = 5000
$ws.Range("L51").Value = 5000
$ws.Range("N51").Value = -5968

$ws.Range("H62").Value = 55561390
$ws.Range("I62").Value = 333333340
$ws.Range("J62").Value = 7001.2
$ws.Range("K62").Value = 333333340
$ws.Range("L62").Value = 7001.2
$ws.Range("M62").Value = -333332716
$ws.Range("N62").Value = -8249.200000000001

$ws.Range("H64").Value = 9771.214
$ws.Range("I64").Value = 7448.5
$ws.Range("K64").Value = 7448.5
$ws.Range("M64").Value = -7200.5

$ws.Range("H65").Value = 55561390
$ws.Range("I65").Value = 333333340
$ws.Range("J65").Value = 7001.2
$ws.Range("K65").Value = 1666666700
$ws.Range("L65").Value = 35006
$ws.Range("M65").Value = -1666663580
$ws.Range("N65").Value = -41246

$ws.Range("H67").Value = 9771.214
$ws.Range("I67").Value = 7448.5
$ws.Range("K67").Value = 7448.5
$ws.Range("M67").Value = -6590.5

$ws.Range("H70").Value = 87875.836
$ws.Range("J70").Value = 5739
$ws.Range("L70").Value = 17217
$ws.Range("N70").Value = -17757

$ws.Range("H73").Value = 87875.836
$ws.Range("J73").Value = 5739
$ws.Range("L73").Value = 17217
$ws.Range("N73").Value = -19089

$ws.Range("H86").Value = 3887.4375
$ws.Range("I86").Value = 4400
$ws.Range("K86").Value = 4400
$ws.Range("M86").Value = -3277

$ws.Range("H89").Value = 3887.4375
$ws.Range("I89").Value = 4400
$ws.Range("K89").Value = 22000
$ws.Range("M89").Value = -16384

$ws.Range("H138").Value = 5025.8
$ws.Range("I138").Value = 3095.077
$ws.Range("J138").Value = 6502.2354
$ws.Range("K138").Value = 9285.231
$ws.Range("L138").Value = 19506.7062
$ws.Range("M138").Value = -4145.231
$ws.Range("N138").Value = -29786.7062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12051241
$ws.Range("I32").Value = 13160092
$ws.Range("K32").Value = 13160092
$ws.Range("M32").Value = -13159805

$ws.Range("H33").Value = 19498.334
$ws.Range("I33").Value = 18495
$ws.Range("J33").Value = 20000
$ws.Range("K33").Value = 18495
$ws.Range("L33").Value = 20000
$ws.Range("M33").Value = -18166
$ws.Range("N33").Value = -20658

$ws.Range("H45").Value = 3281.7058
$ws.Range("I45").Value = 1106.4546
$ws.Range("K45").Value = 1106.4546
$ws.Range("M45").Value = -729.4546

$ws.Range("H53").Value = 46677.668
$ws.Range("J53").Value = 46677.668
$ws.Range("L53").Value = 46677.668
$ws.Range("N53").Value = -48041.668

$ws.Range("H61").Value = 6717.2144
$ws.Range("I61").Value = 6717.2144
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6717.2144
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -6505.2144
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 13890843
$ws.Range("I74").Value = 14494488
$ws.Range("J74").Value = 7000
$ws.Range("K74").Value = 14494488
$ws.Range("L74").Value = 7000
$ws.Range("M74").Value = -14493614
$ws.Range("N74").Value = -8748

$ws.Range("H77").Value = 13890843
$ws.Range("I77").Value = 14494488
$ws.Range("J77").Value = 7000
$ws.Range("K77").Value = 72472440
$ws.Range("L77").Value = 35000
$ws.Range("M77").Value = -72468072
$ws.Range("N77").Value = -43736

$ws.Range("H122").Value = 1812.2424
$ws.Range("I122").Value = 1848.1786
$ws.Range("J122").Value = 1611
$ws.Range("K122").Value = 5544.5358
$ws.Range("L122").Value = 4833
$ws.Range("M122").Value = -3094.5358
$ws.Range("N122").Value = -9733

$ws.Range("H132").Value = 3999.8845
$ws.Range("I132").Value = 1449.8182
$ws.Range("J132").Value = 18025.25
$ws.Range("K132").Value = 4349.4546
$ws.Range("L132").Value = 54075.75
$ws.Range("M132").Value = -1819.4546
$ws.Range("N132").Value = -59135.75

$ws.Range("H136").Value = 6717.2144
$ws.Range("I136").Value = 6717.2144
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 20151.6432
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -17601.6432
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1229
$ws.Range("I94").Value = 778.9375
$ws.Range("J94").Value = 3629.3333
$ws.Range("K94").Value = 778.9375
$ws.Range("L94").Value = 3629.3333
$ws.Range("M94").Value = -327.9375
$ws.Range("N94").Value = -4531.3333

$ws.Range("H107").Value = 1731.9565
$ws.Range("I107").Value = 1492.5
$ws.Range("J107").Value = 7000
$ws.Range("K107").Value = 1492.5
$ws.Range("L107").Value = 7000
$ws.Range("M107").Value = 427.5
$ws.Range("N107").Value = -10840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3298.0588
$ws.Range("I16").Value = 2211.0908
$ws.Range("K16").Value = 2211.0908
$ws.Range("M16").Value = -1924.0908

$ws.Range("H31").Value = 41534.703
$ws.Range("I31").Value = 1597.6666
$ws.Range("K31").Value = 1597.6666
$ws.Range("M31").Value = -1302.6666

$ws.Range("H34").Value = 41534.703
$ws.Range("I34").Value = 1597.6666
$ws.Range("K34").Value = 1597.6666
$ws.Range("M34").Value = -1395.6666

$ws.Range("H62").Value = 10974
$ws.Range("I62").Value = 4513.5
$ws.Range("J62").Value = 13127.5
$ws.Range("K62").Value = 4513.5
$ws.Range("L62").Value = 13127.5
$ws.Range("M62").Value = -3889.5
$ws.Range("N62").Value = -14375.5

$ws.Range("H65").Value = 10974
$ws.Range("I65").Value = 4513.5
$ws.Range("J65").Value = 13127.5
$ws.Range("K65").Value = 22567.5
$ws.Range("L65").Value = 65637.5
$ws.Range("M65").Value = -19447.5
$ws.Range("N65").Value = -71877.5

$ws.Range("H107").Value = 770.5
$ws.Range("I107").Value = 554.75
$ws.Range("K107").Value = 554.75
$ws.Range("M107").Value = 1365.25

$ws.Range("H113").Value = 3298.0588
$ws.Range("I113").Value = 2211.0908
$ws.Range("K113").Value = 2211.0908
$ws.Range("M113").Value = -41.09079999999994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 233333.92
$ws.Range("J37").Value = 233333.92
$ws.Range("L37").Value = 700001.76
$ws.Range("N37").Value = -700225.76

$ws.Range("H129").Value = 3790450.8
$ws.Range("I129").Value = 778.0909
$ws.Range("J129").Value = 7580123.5
$ws.Range("K129").Value = 2334.2727
$ws.Range("L129").Value = 22740370.5
$ws.Range("M129").Value = 2665.7273
$ws.Range("N129").Value = -22750370.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9211.666999999999
$ws.Range("I80").Value = 6333.3335
$ws.Range("K80").Value = 6333.3335
$ws.Range("M80").Value = -5335.3335

$ws.Range("H83").Value = 9211.666999999999
$ws.Range("I83").Value = 6333.3335
$ws.Range("K83").Value = 31666.6675
$ws.Range("M83").Value = -26674.6675

$ws.Range("H126").Value = 4637.467
$ws.Range("I126").Value = 3154.9
$ws.Range("K126").Value = 9464.700000000001
$ws.Range("M126").Value = -6994.700000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3579.9
$ws.Range("J46").Value = 4349.875
$ws.Range("L46").Value = 4349.875
$ws.Range("N46").Value = -4725.875

$ws.Range("H55").Value = 2377.3125
$ws.Range("I55").Value = 305.66666
$ws.Range("J55").Value = 3620.3
$ws.Range("K55").Value = 305.66666
$ws.Range("L55").Value = 3620.3
$ws.Range("M55").Value = -132.66666
$ws.Range("N55").Value = -3966.3

$ws.Range("H136").Value = 6199.7144
$ws.Range("I136").Value = 3066.3333
$ws.Range("J136").Value = 25000
$ws.Range("K136").Value = 9198.999899999999
$ws.Range("L136").Value = 75000
$ws.Range("M136").Value = -6648.999899999999
$ws.Range("N136").Value = -80100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 24500

$ws.Range("H62").Value = 7900
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 7900
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H126").Value = 1984.7576
$ws.Range("I126").Value = 1816.5
$ws.Range("J126").Value = 2243.6155
$ws.Range("K126").Value = 5449.5
$ws.Range("L126").Value = 6730.8465
$ws.Range("M126").Value = -2979.5
$ws.Range("N126").Value = -11670.8465

$ws.Range("H132").Value = 8224.585999999999
$ws.Range("I132").Value = 3715.25
$ws.Range("K132").Value = 11145.75
$ws.Range("M132").Value = -8615.75

$ws.Range("H136").Value = 2016.45
$ws.Range("I136").Value = 2016.45
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6049.35
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3499.35
$ws.Range("N136").ClearContents()
